$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The oldest year (row 2, "2009年") drops out of the series; every
# remaining year shifts up one row (2010..2020 -> rows 2..12), and a new
# year ("2021年") is appended as the new last row (13).
$ws.Rows(2).Delete()

# Give the new row's label cell the same formatting (bold/border/centered)
# used by the rest of the year column before filling in its value.
$ws.Range("A12").Copy($ws.Range("A13"))

# Populate the new last row with the 2021 figures.
$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 105.4672
$ws.Range("E13").Value = 243.4498
$ws.Range("G13").Value = 699.2403
$ws.Range("I13").Value = 135.6675
$ws.Range("J13").Value = 168.9282
$ws.Range("O13").Value = 274.3954
